# Update "想去人数" (want-to-go count) figures in column F.
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both list the same
# events, so the same counts must be bumped on both sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# "展览" sheet: row -> new F value
$wsExhibit.Range("F2").Value  = 454
$wsExhibit.Range("F4").Value  = 74
$wsExhibit.Range("F5").Value  = 5161
$wsExhibit.Range("F6").Value  = 175
$wsExhibit.Range("F7").Value  = 57
$wsExhibit.Range("F8").Value  = 97
$wsExhibit.Range("F9").Value  = 336
$wsExhibit.Range("F10").Value = 3
$wsExhibit.Range("F11").Value = 61

# "全部类型" sheet: same events, different row offsets
$wsAll.Range("F2").Value  = 454
$wsAll.Range("F8").Value  = 74
$wsAll.Range("F9").Value  = 5161
$wsAll.Range("F10").Value = 175
$wsAll.Range("F11").Value = 57
$wsAll.Range("F12").Value = 97
$wsAll.Range("F14").Value = 336
$wsAll.Range("F15").Value = 3
$wsAll.Range("F16").Value = 61
